$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) First paragraph: append "  (This is a change - Version for branch
#    alternate)" after the existing sentence. The trailing text is typed in
#    three separate red (C00000) bursts, matching how Word keeps each
#    "typing session" as its own run.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$r = $p1.Range
$r.End = $r.End - 1                 # exclude the paragraph mark
$pos = $r.End

$r.InsertAfter("  ")                # two trailing spaces, default formatting
$pos = $pos + 2

$enDash = [char]0x2013
$redColor = 192                     # RGB(0xC0,0x00,0x00) -> wdColor for C00000

$chunk1 = "(This is a change " + $enDash + " Ve"
$rr = $d.Range($pos, $pos)
$rr.InsertAfter($chunk1)
$d.Range($pos, $pos + $chunk1.Length).Font.Color = $redColor
$pos = $pos + $chunk1.Length

$chunk2 = "rsion for branch alternate"
$rr = $d.Range($pos, $pos)
$rr.InsertAfter($chunk2)
$d.Range($pos, $pos + $chunk2.Length).Font.Color = $redColor
$pos = $pos + $chunk2.Length

$chunk3 = ")"
$rr = $d.Range($pos, $pos)
$rr.InsertAfter($chunk3)
$d.Range($pos, $pos + $chunk3.Length).Font.Color = $redColor
$pos = $pos + $chunk3.Length

# ---------------------------------------------------------------------------
# 2) Append a new, empty, shaded paragraph after the final paragraph of the
#    document (before the sectPr).
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$tail = $lastPara.Range
$tail.Collapse(0)                   # collapse to the end of the last paragraph
$tail.Style = "Normal"
$tail.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Style = $d.Styles(-1)      # wdStyleNormal, applied without stamping an explicit pStyle
$newPara.Shading.Texture = 0
$newPara.Shading.ForegroundPatternColor = -16777216
$newPara.Shading.BackgroundPatternColor = 0xF9F9F9

Write-Output "done"
